# "Changed USB to differential pair" - add a sample-size worksheet below the
# existing BRG/baud table: bits/channel * channels/card * cards = bits/sample,
# then bits/sample / 8 = bytes/sample.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the label/value columns so the new rows are readable.
$ws.Columns.Item(1).ColumnWidth = 17.5703125
$ws.Columns.Item(2).ColumnWidth = 12.7109375

# Row 12: Bits per channel
$ws.Range("A12").Value = "Bits per channel"
$ws.Range("B12").Value = 16

# Row 13: Channels per card
$ws.Range("A13").Value = "Channels per card"
$ws.Range("B13").Value = 2

# Row 14: Cards
$ws.Range("A14").Value = "Cards"
$ws.Range("B14").Value = 12

# Row 15: Bits per sample
$ws.Range("A15").Value = "Bits per sample"
$ws.Range("B15").Formula = "=B12*B13*B14"

# Row 16: Bytes per sample
$ws.Range("A16").Value = "Bytes per sample"
$ws.Range("B16").Formula = "=B15/8"

# Update selection to match diff
$ws.Range("B20").Select()
